$d = $word.ActiveDocument
$d.Content.Find.Execute("Max Musterman", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Max Mustermann", 2)
